$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the account-number column is treated as text so the long
# all-digit value is not coerced into a number.
$ws.Range("C2:C3").NumberFormat = "@"

# Row 2: first virement entry for BAKKALI MOHAMED
$ws.Range("A2").Value = "BAKKALI MOHAMED"
$ws.Range("B2").Value = "B12346"
$ws.Range("C2").Value = "78017053636372722873881919"
$ws.Range("D2").Value = "HASSAN 2"
$ws.Range("E2").Value = "CIH"
$ws.Range("F2").Value = "Direction régionale"
$ws.Range("G2").Value = "988/DIRECTION CAPITAL SOFT"
$ws.Range("H2").Value = "mensuelle"
$ws.Range("I2").Value = 40000
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 38000

# Row 3: second virement entry for BAKKALI MOHAMED
$ws.Range("A3").Value = "BAKKALI MOHAMED"
$ws.Range("B3").Value = "B12346"
$ws.Range("C3").Value = "78017053636372722873881919"
$ws.Range("D3").Value = "HASSAN 2"
$ws.Range("E3").Value = "CIH"
$ws.Range("F3").Value = "Direction régionale"
$ws.Range("G3").Value = "988/DIRECTION CAPITAL SOFT"
$ws.Range("H3").Value = "mensuelle"
$ws.Range("I3").Value = 88000
$ws.Range("J3").Value = 400
$ws.Range("K3").Value = 87600

# Row 4: totals row
$ws.Range("A4").Value = " "
$ws.Range("B4").Value = " "
$ws.Range("C4").Value = " "
$ws.Range("D4").Value = " "
$ws.Range("E4").Value = " "
$ws.Range("F4").Value = " "
$ws.Range("G4").Value = " "
$ws.Range("H4").Value = " "
$ws.Range("I4").Value = 128000
$ws.Range("J4").Value = 2400
$ws.Range("K4").Value = 125600
